# Optimize the display of the results
# - Append three new weekly performance rows (6, 7, 8) to the portfolio
#   performance table on Sheet1.
# - Keep the date column formatted the same way as the existing date cells.
# - Widen column A so the date values are fully visible.
# - Leave the selection on the sheet where the author left off (E11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: Date, S&P500, ChatGPT, Gemini (portfolio values)
$newRows = @(
    @{ Row = 6; Date = 45808; SP500 = 5911; ChatGPT = 1057899; Gemini = 1105496 },
    @{ Row = 7; Date = 45815; SP500 = 6000; ChatGPT = 1087661; Gemini = 1124766 },
    @{ Row = 8; Date = 45822; SP500 = 5976; ChatGPT = 1079914; Gemini = 1094591 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.SP500
    $ws.Cells.Item($row, 3).Value = $r.ChatGPT
    $ws.Cells.Item($row, 4).Value = $r.Gemini

    # Match the date number format already used by the rows above (A2:A5)
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
}

# Widen the Date column so the new values display cleanly
$ws.Columns.Item(1).ColumnWidth = 18.83

# Restore the sheet's last-used selection
$ws.Range("E11").Select() | Out-Null
